$wb = $excel.ActiveWorkbook

# --- Sheet: Summary (5 cell updates) ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = 195407.0288009005
$ws.Range("B7").Value = 9654200.856693137
$ws.Range("B8").Value = 17785260.94054101
$ws.Range("B10").Value = 5582999.386126758

# --- Sheet: Costs and Revenues (37 cell updates) ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 512047.7116981491
$ws.Range("C2").Value = 512047.7116981489
$ws.Range("E2").Value = 142900.6851873295
$ws.Range("G2").Value = 142900.6851873295
$ws.Range("H2").Value = 142900.6851873295
$ws.Range("I2").Value = 142900.6851873295
$ws.Range("J2").Value = 142900.6851873295
$ws.Range("K2").Value = 142900.6851873295
$ws.Range("M2").Value = 142900.6851873295
$ws.Range("N2").Value = 142900.6851873295
$ws.Range("O2").Value = 142900.6851873295
$ws.Range("P2").Value = 142900.6851873295
$ws.Range("E4").Value = 79636.80000449967
$ws.Range("G4").Value = 79636.80000449967
$ws.Range("H4").Value = 79636.80000449967
$ws.Range("I4").Value = 79636.80000449967
$ws.Range("J4").Value = 79636.80000449967
$ws.Range("K4").Value = 79636.80000449967
$ws.Range("M4").Value = 79636.80000449967
$ws.Range("N4").Value = 79636.80000449967
$ws.Range("O4").Value = 79636.80000449967
$ws.Range("P4").Value = 79636.80000449967
$ws.Range("B6").Value = 10760.01598773675
$ws.Range("C6").Value = 10760.01598773658
$ws.Range("D6").Value = 10760.01598773669
$ws.Range("E6").Value = -88257.74376129247
$ws.Range("F6").Value = 44842.25623870757
$ws.Range("G6").Value = 44842.25623870757
$ws.Range("H6").Value = 44842.25623870757
$ws.Range("I6").Value = 44842.25623870757
$ws.Range("J6").Value = 44842.25623870757
$ws.Range("K6").Value = 44842.25623870757
$ws.Range("L6").Value = 44842.25623870757
$ws.Range("M6").Value = 44842.25623870757
$ws.Range("N6").Value = 44842.25623870757
$ws.Range("O6").Value = 44842.25623870757
$ws.Range("P6").Value = 44842.25623870757

# --- Sheet: Fed-in Capacity (54 cell updates) ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("J13").Value = 33.63624132272333
$ws.Range("K13").Value = 94.30397654773019
$ws.Range("L13").Value = 90.4687457914608
$ws.Range("M13").Value = 92.09541281912071
$ws.Range("N13").Value = 81.96869489115805
$ws.Range("O13").Value = 96.22962838366004
$ws.Range("P13").Value = 101.5955875616828
$ws.Range("Q13").Value = 65.34295837775146
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("J20").Value = 124.5190384721106
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("O21").Value = 57.81213424001893
$ws.Range("P21").Value = 65.92768427608706
$ws.Range("L22").Value = 90.4687457914608
$ws.Range("N23").Value = 110.5750244233121
$ws.Range("L26").Value = 130.6648563030561
$ws.Range("M26").Value = 113.4004983079896
$ws.Range("N26").Value = 110.5750244233121
$ws.Range("O26").Value = 117.8828208804077
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("L27").Value = 61.18167021676314
$ws.Range("M27").Value = 51.84373129681028
$ws.Range("N27").Value = 38.66169381481656
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("L28").Value = 90.4687457914608
$ws.Range("L29").Value = 130.6648563030561
$ws.Range("N29").Value = 110.5750244233121
$ws.Range("N30").Value = 38.66169381481656
$ws.Range("O30").Value = 57.81213424001893
$ws.Range("P30").Value = 65.92768427608706
$ws.Range("M36").Value = 51.84373129681028
$ws.Range("N36").Value = 38.66169381481656
$ws.Range("J38").Value = 124.5190384721106
$ws.Range("N39").Value = 38.66169381481656
$ws.Range("L41").Value = 130.6648563030561
$ws.Range("M41").Value = 113.4004983079896
$ws.Range("N41").Value = 110.5750244233121
$ws.Range("O41").Value = 117.8828208804077
$ws.Range("L42").Value = 61.18167021676314
$ws.Range("M42").Value = 51.84373129681028
$ws.Range("N42").Value = 38.66169381481656
$ws.Range("O42").Value = 57.81213424001893
$ws.Range("J44").Value = 124.5190384721106
$ws.Range("K44").Value = 135.370731907559
$ws.Range("L44").Value = 130.6648563030561
$ws.Range("M44").Value = 113.4004983079896
$ws.Range("N44").Value = 110.5750244233121
$ws.Range("O44").Value = 117.8828208804077
$ws.Range("P44").Value = 135.4597561231036
$ws.Range("Q44").Value = 150.3839754851235
$ws.Range("M45").Value = 51.84373129681028
$ws.Range("O45").Value = 57.81213424001893

# --- Sheet: Unmet Demand (54 cell updates) ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("J13").Value = 72.23757736389061
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 61.14583096471014
$ws.Range("P18").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("O45").Value = 0

# --- Sheet: Household Surplus (10 cell updates) ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 336858.0333416928
$ws.Range("B7").Value = 336858.0333416928
$ws.Range("B8").Value = 336858.0333416928
$ws.Range("B9").Value = 336858.0333416928
$ws.Range("B10").Value = 336858.0333416928
$ws.Range("B11").Value = 336858.0333416928
$ws.Range("B13").Value = 336858.0333416928
$ws.Range("B14").Value = 336858.0333416928
$ws.Range("B15").Value = 336858.0333416928
$ws.Range("B16").Value = 336858.0333416928

Write-Output "Applied all cell updates"